$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "28.050.96"
$ws.Range("E2").Value = "  +3.27%  "
$ws.Range("D3").Value = "1.574.37"
$ws.Range("E3").Value = "  +0.29%  "
$ws.Range("E4").Value = "  -1.09%  "
$ws.Range("D5").Value = "'212.77"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +0.72%  "
$ws.Range("E6").Value = "  -0.09%  "
$ws.Range("D7").Value = "'0.998"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -1.11%  "
$ws.Range("D8").Value = "'23.17"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +5.24%  "
$ws.Range("D9").Value = "'0.251"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +0.67%  "
$ws.Range("D10").Value = "'0.0598"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -0.26%  "
$ws.Range("D11").Value = "'0.0883"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +1.84%  "
$ws.Range("D12").Value = "1.798.89"
$ws.Range("E12").Value = "  +0.30%  "
$ws.Range("D13").Value = "1.574.51"
$ws.Range("E13").Value = "  +0.37%  "
$ws.Range("E14").Value = "  -0.63%  "
$ws.Range("D16").Value = "28.029.26"
$ws.Range("E16").Value = "  +3.31%  "
$ws.Range("D17").Value = "'63.51"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +2.07%  "
$ws.Range("D18").Value = "'228.75"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +5.97%  "
$ws.Range("D19").Value = "0.0₃0706"
$ws.Range("E19").Value = "  +0.36%  "
$ws.Range("D20").Value = "'7.45"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +0.64%  "
$ws.Range("D21").Value = "'0.998"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -1.16%  "
$ws.Range("D22").Value = "'4.12"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -0.75%  "
$ws.Range("D23").Value = "'9.33"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +1.32%  "
$ws.Range("D24").Value = "'1.94"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +0.30%  "
$ws.Range("D25").Value = "'152.40"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -0.84%  "
$ws.Range("D26").Value = "'15.24"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +1.21%  "
$ws.Range("E27").Value = "  -0.91%  "
$ws.Range("E28").Value = "  +0.17%  "
$ws.Range("D29").Value = "'0.999"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -1.08%  "
$ws.Range("E30").Value = "  +0.21%  "
$ws.Range("E31").Value = "  +0.06%  "
$ws.Range("E32").Value = "  -0.49%  "
$ws.Range("E33").Value = "  -1.95%  "
$ws.Range("D34").Value = "1.417.09"
$ws.Range("E34").Value = "  -2.44%  "
$ws.Range("E35").Value = "  -1.47%  "
$ws.Range("E36").Value = "  -4.69%  "
$ws.Range("E37").Value = "  -1.39%  "
$ws.Range("E38").Value = "  -0.11%  "
$ws.Range("D39").Value = "'0.540"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +1.21%  "
$ws.Range("D40").Value = "'2.47"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +4.22%  "
$ws.Range("E41").Value = "  -0.40%  "
$ws.Range("E42").Value = "  -1.21%  "
$ws.Range("E43").Value = "  -2.76%  "
$ws.Range("D44").Value = "'0.974"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -2.53%  "
$ws.Range("D45").Value = "'1.81"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +4.88%  "
$ws.Range("E46").Value = "  -1.35%  "
$ws.Range("D47").Value = "1.712.03"
$ws.Range("E47").Value = "  +0.61%  "
$ws.Range("D48").Value = "'86.90"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +0.96%  "
$ws.Range("E49").Value = "  +3.10%  "
$ws.Range("E50").Value = "  +0.71%  "
$ws.Range("E51").Value = "  -1.85%  "
